$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value ($null means clear the cell)
$changes = @{
    "A1" = $null; "C1" = $null; "F1" = 9; "H1" = $null; "I1" = $null;
    "A2" = 8; "C2" = $null; "D2" = 6; "E2" = 3; "G2" = $null; "I2" = 4;
    "B3" = 3; "D3" = $null; "E3" = $null; "F3" = 8; "G3" = 6; "H3" = 9; "I3" = 1;
    "A4" = $null; "C4" = $null; "D4" = 7; "E4" = $null; "F4" = 5; "G4" = 3; "H4" = $null; "I4" = 9;
    "B5" = 8; "H5" = $null; "I5" = 5;
    "B6" = 7; "C6" = 5; "E6" = 9; "F6" = $null; "G6" = 1;
    "A7" = 1; "B7" = 4; "C7" = 8; "D7" = $null; "E7" = 5; "F7" = $null; "G7" = $null;
    "B8" = 6; "F8" = 4; "G8" = 5;
    "A9" = 5; "B9" = $null; "C9" = $null; "E9" = 1; "G9" = $null; "I9" = 6;
}

foreach ($addr in $changes.Keys) {
    $val = $changes[$addr]
    if ($null -eq $val) {
        $ws.Range($addr).ClearContents()
    } else {
        $ws.Range($addr).Value = $val
    }
}
